$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for the rows that changed.
$values = @{
    2  = -3
    4  = -1
    5  = -2
    6  = 1
    8  = -3
    9  = -4
    10 = 5
    11 = -2
    13 = -5
    14 = 1
    15 = -3
    16 = 1
    17 = 4
    18 = -5
    19 = 6
    20 = 5
    21 = -4
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
